$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 08:18:19"
$ws.Range("E3").Value = "2026-02-20 08:18:21"
$ws.Range("E4").Value = "2026-02-20 08:18:24"
$ws.Range("E5").Value = "2026-02-20 08:18:26"
$ws.Range("E6").Value = "2026-02-20 08:18:29"
$ws.Range("E7").Value = "2026-02-20 08:18:31"
$ws.Range("E8").Value = "2026-02-20 08:18:34"
$ws.Range("E9").Value = "2026-02-20 08:18:36"
$ws.Range("E10").Value = "2026-02-20 08:18:39"
$ws.Range("E11").Value = "2026-02-20 08:18:41"
$ws.Range("E12").Value = "2026-02-20 08:18:43"
$ws.Range("E13").Value = "2026-02-20 08:18:46"
$ws.Range("E14").Value = "2026-02-20 08:18:48"
$ws.Range("E15").Value = "2026-02-20 08:18:51"
$ws.Range("E16").Value = "2026-02-20 08:18:53"
$ws.Range("E17").Value = "2026-02-20 08:18:56"
$ws.Range("E18").Value = "2026-02-20 08:18:58"
$ws.Range("E19").Value = "2026-02-20 08:19:00"
$ws.Range("E20").Value = "2026-02-20 08:19:03"
$ws.Range("E21").Value = "2026-02-20 08:19:05"
$ws.Range("E22").Value = "2026-02-20 08:19:08"
$ws.Range("E23").Value = "2026-02-20 08:19:10"
$ws.Range("E24").Value = "2026-02-20 08:19:13"
$ws.Range("E25").Value = "2026-02-20 08:19:15"
$ws.Range("E26").Value = "2026-02-20 08:19:17"
$ws.Range("E27").Value = "2026-02-20 08:19:20"
$ws.Range("E28").Value = "2026-02-20 08:19:22"
$ws.Range("E29").Value = "2026-02-20 08:19:25"
$ws.Range("E30").Value = "2026-02-20 08:19:27"
$ws.Range("E31").Value = "2026-02-20 08:19:29"
$ws.Range("E32").Value = "2026-02-20 08:19:32"
$ws.Range("E33").Value = "2026-02-20 08:19:34"
$ws.Range("E34").Value = "2026-02-20 08:19:37"
$ws.Range("E35").Value = "2026-02-20 08:19:39"
$ws.Range("E36").Value = "2026-02-20 08:19:42"
$ws.Range("E37").Value = "2026-02-20 08:19:44"
$ws.Range("E38").Value = "2026-02-20 08:19:47"
$ws.Range("E39").Value = "2026-02-20 08:19:49"
$ws.Range("E40").Value = "2026-02-20 08:19:52"
$ws.Range("E41").Value = "2026-02-20 08:19:54"
$ws.Range("E42").Value = "2026-02-20 08:19:56"
$ws.Range("E43").Value = "2026-02-20 08:19:59"
$ws.Range("E44").Value = "2026-02-20 08:20:01"
$ws.Range("E45").Value = "2026-02-20 08:20:04"
$ws.Range("E46").Value = "2026-02-20 08:20:06"
